$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.719.57'
$ws.Range("E2").Value = '  +3.18%  '
$ws.Range("D3").Value = '2.552.13'
$ws.Range("E3").Value = '  +3.63%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '580.24'
$ws.Range("E5").Value = '  +1.23%  '
$ws.Range("D6").Value = '152.82'
$ws.Range("E6").Value = '  +3.51%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +1.00%  '
$ws.Range("D9").Value = '2.552.94'
$ws.Range("E9").Value = '  +3.57%  '
$ws.Range("E10").Value = '  +0.99%  '
$ws.Range("E11").Value = '  -1.70%  '
$ws.Range("E12").Value = '  -0.40%  '
$ws.Range("D13").Value = '0.354'
$ws.Range("E13").Value = '  +0.18%  '
$ws.Range("D14").Value = '29.12'
$ws.Range("E14").Value = '  +0.15%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '3.016.42'
$ws.Range("E15").Value = '  +3.55%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000179'
$ws.Range("E16").Value = '  +1.63%  '
$ws.Range("D17").Value = '64.587.26'
$ws.Range("E17").Value = '  +3.12%  '
$ws.Range("D18").Value = '2.552.92'
$ws.Range("E18").Value = '  +3.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.00'
$ws.Range("E19").Value = '  +1.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.00'
$ws.Range("E20").Value = '  +0.92%  '
$ws.Range("D21").Value = '340.89'
$ws.Range("E21").Value = '  +4.60%  '
$ws.Range("D22").Value = '4.28'
$ws.Range("E22").Value = '  +3.58%  '
$ws.Range("E23").Value = '  +2.66%  '
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").Value = '10.08'
$ws.Range("E25").Value = '  +1.49%  '
$ws.Range("D26").Value = '65.96'
$ws.Range("D27").Value = '626.56'
$ws.Range("E27").Value = '  -2.26%  '
$ws.Range("E28").Value = '  +6.44%  '
$ws.Range("D29").Value = '2.716.31'
$ws.Range("E29").Value = '  +5.04%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  +0.13%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").Value = '1.48'
$ws.Range("E31").Value = '  +3.38%  '
$ws.Range("E32").Value = '  +2.09%  '
$ws.Range("D33").Value = '1.87'
$ws.Range("E33").Value = '  +2.66%  '
$ws.Range("E34").Value = '  +4.09%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D37").Value = '4.85'
$ws.Range("E37").Value = '  +2.30%  '
$ws.Range("D38").Value = '5.58'
$ws.Range("E38").Value = '  +5.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '155.70'
$ws.Range("E39").Value = '  +3.26%  '
$ws.Range("D40").Value = '19.01'
$ws.Range("E40").Value = '  +2.33%  '
$ws.Range("E41").Value = '  +4.41%  '
$ws.Range("D42").Value = '0.372'
$ws.Range("E42").Value = '  +1.25%  '
$ws.Range("E43").Value = '  +4.27%  '
$ws.Range("D44").Value = '42.07'
$ws.Range("E44").Value = '  +0.74%  '
$ws.Range("D45").Value = '161.57'
$ws.Range("E45").Value = '  +5.60%  '
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("E47").Value = '  -3.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.70'
$ws.Range("E49").Value = '  +2.82%  '
$ws.Range("D50").Value = '21.54'
$ws.Range("E50").Value = '  +5.95%  '
$ws.Range("D51").Value = '0.628'
$ws.Range("E51").Value = '  +3.61%  '
